$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update B2/C2, and fill in D2/E2/F2 (previously empty) ---
$ws.Range("B2").Value = "NSE:360ONE"
$ws.Range("C2").Value = "NSE:ABB"
$ws.Range("D2").Value = "NSE:MUTHOOTFIN"
$ws.Range("E2").Value = "NSE:APLAPOLLO"
$ws.Range("F2").Value = "NSE:BOSCHLTD"

# --- Rows 3-21: update column B and C values ---
$ws.Range("B3").Value = "NSE:ALPHAGEO"
$ws.Range("C3").Value = "NSE:ADANIPORTS"

$ws.Range("B4").Value = "NSE:APOLSINHOT"
$ws.Range("C4").Value = "NSE:ALMONDZ"

$ws.Range("B5").Value = "NSE:ARMANFIN"
$ws.Range("C5").Value = "NSE:BANKA"

$ws.Range("B6").Value = "NSE:BALMLAWRIE"
$ws.Range("C6").Value = "NSE:BCLIND"

$ws.Range("B7").Value = "NSE:BHAGCHEM"
$ws.Range("C7").Value = "NSE:CENTEXT"

$ws.Range("B8").Value = "NSE:BOSCHLTD"
$ws.Range("C8").Value = "NSE:CENTRALBK"

$ws.Range("B9").Value = "NSE:BTML"
$ws.Range("C9").Value = "NSE:COSMOFIRST"

$ws.Range("B10").Value = "NSE:CELLO"
$ws.Range("C10").Value = "NSE:DANGEE"

$ws.Range("B11").Value = "NSE:ESABINDIA"
$ws.Range("C11").Value = "NSE:E2E"

$ws.Range("B12").Value = "NSE:GMBREW"
$ws.Range("C12").Value = "NSE:EUROTEXIND"

$ws.Range("B13").Value = "NSE:GMMPFAUDLR"
$ws.Range("C13").Value = "NSE:GMDCLTD"

$ws.Range("B14").Value = "NSE:METROPOLIS"
$ws.Range("C14").Value = "NSE:GNFC"

$ws.Range("B15").Value = "NSE:MRPL"
$ws.Range("C15").Value = "NSE:GREENPOWER"

$ws.Range("B16").Value = "NSE:NIPPOBATRY"
$ws.Range("C16").Value = "NSE:IDBI"

$ws.Range("B17").Value = "NSE:ORIENTCER"
$ws.Range("C17").Value = "NSE:IFCI"

$ws.Range("B18").Value = "NSE:PARSVNATH"
$ws.Range("C18").Value = "NSE:JINDALSAW"

$ws.Range("B19").Value = "NSE:PDMJEPAPER"
$ws.Range("C19").Value = "NSE:JIOFIN"

$ws.Range("B20").Value = "NSE:PDSL"
$ws.Range("C20").Value = "NSE:KTKBANK"

$ws.Range("B21").Value = "NSE:RPSGVENT"
$ws.Range("C21").Value = "NSE:NCC"

# --- Row 22: B22 cleared (no longer has a ticker), C22 now populated ---
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "NSE:NDRAUTO"

# --- New rows 23-25: copy the row-index cell formatting from A22, then fill values ---
$ws.Range("A22").Copy()
$ws.Range("A23:A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A23").Value = 21
$ws.Range("C23").Value = "NSE:PERSISTENT"

$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "NSE:PSPPROJECT"

$ws.Range("A25").Value = 23
$ws.Range("C25").Value = "NSE:RAMASTEEL"
